$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1
$ws.Range("B1").Value = "Untitled 1"
$ws.Range("C1").Value = "Untitled 2"
$ws.Range("D1").Value = "Untitled 3"

# Update A2 and add B2:D2, matching A2's existing style
$ws.Range("A2").Value = 0.000000
$ws.Range("B2").Value = 0.000007
$ws.Range("C2").Value = -0.656128
$ws.Range("D2").Value = 0.656128

$ws.Range("B2:D2").HorizontalAlignment = $ws.Range("A2").HorizontalAlignment
